# Generate Report for Handback
# Update the Handoff/Handback datetimes for the first file row
# (3af81ce8-...) on the "zh-cn" and "de-de" sheets to reflect a
# re-run of the handback report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 and 3 both show the same handoff/handback times
# for the 3af81ce8 entry.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2:E3").Value = "2016-03-24 18:26:25"
$wsZh.Range("H2:H3").Value = "2016-03-24 18:27:02"

# de-de sheet: same pair of cells, different timestamps.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2:E3").Value = "2016-03-24 18:26:30"
$wsDe.Range("H2:H3").Value = "2016-03-24 18:27:13"
